$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 4 ("Doktoratsheuriger"): shape "CustomShape 2" (the bulleted body copy)
# ---------------------------------------------------------------------------
$s  = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# 1) Remove the whole "Online teaching coffee in the making" bullet paragraph.
$bullet = $tr.Paragraphs(3, 1)
$bullet.Delete()

# 2) "TBD (Fall)" -> "TBD (~End of June)" (after the delete above this is now
#    paragraph 5, was paragraph 6).
$tbdPara = $tr.Paragraphs(5, 1)

# Replace the whole visible run text in one shot so the new text keeps the
# original run's formatting (lang="de-AT" sz=2800 Calibri/DejaVu Sans etc.)
# instead of picking up blank/default formatting.
$whole = $tr.Characters($tbdPara.Start, 10)
$whole.Text = "TBD (~End of June)"

# Re-select just "of" inside the freshly written text and re-assign it so it
# becomes its own run, split off from the surrounding text.
$tbdPara2 = $tr.Paragraphs(5, 1)
$ofRange = $tr.Characters($tbdPara2.Start + 10, 2)
$ofRange.Text = "of"
